$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the data (A2:C26) by the Model column (C) only, instead of by
# DietSource (A) then Model (C) -- drops the old A2:A26 sort condition.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C26")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:C26"))
$ws.Sort.Header = 1
$ws.Sort.Apply() | Out-Null

# Widen the Estimate column and center its contents (header first, then
# the numeric data, so the header's plain-centered style is interned
# before the number-format+centered style used by the data cells).
$ws.Columns("B").ColumnWidth = 13.67
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B2:B26").NumberFormat = "0.00"
$ws.Range("B2:B26").HorizontalAlignment = -4108

# Move the active selection
$ws.Range("C7").Select() | Out-Null
